$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving a numeric-looking string value must be forced to Text format
# so Excel does not silently convert them to numbers (losing formatting like trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "37.857.21"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "2.033.85"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("D5").Value = "227.61"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("D7").Value = "60.15"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D10").Value = "0.0817"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "2.335.45"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "14.46"
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("D14").Value = "21.19"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "0.759"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").Value = "2.030.26"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "37.783.80"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "69.78"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "5.88"
$ws.Range("E20").Value = "  -6.55%  "
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("D22").Value = "223.80"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "2.26"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "9.36"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "167.32"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("E28").Value = "  -2.91%  "
$ws.Range("D29").Value = "18.85"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("E30").Value = "  -4.17%  "
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("E32").Value = "  +8.72%  "
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "4.49"
$ws.Range("E35").Value = "  -2.47%  "
$ws.Range("D36").Value = "6.36"
$ws.Range("E36").Value = "  +4.18%  "
$ws.Range("D37").Value = "2.28"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("D38").Value = "3.33"
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.533.87"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "17.58"
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "96.12"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").Value = "0.0912"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("E46").Value = "  -3.18%  "
$ws.Range("D47").Value = "3.98"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "7.14"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "2.96"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").Value = "2.224.77"
$ws.Range("E51").Value = "  -1.14%  "
